# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force the cell to be treated as literal text even when the
    # string looks like a number (e.g. "0.999"), then restore the
    # default "Normal" cell style so no stray formatting is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "63.126.64"
$ws.Range("E2").Value = "  +2.86%  "
$ws.Range("D3").Value = "2.954.74"
$ws.Range("E3").Value = "  +0.91%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "594.92"
$ws.Range("E5").Value = "  -0.49%  "
Set-TextValue "D6" "148.67"
$ws.Range("E6").Value = "  +2.60%  "
Set-TextValue "D7" "0.998"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "2.953.07"
$ws.Range("E8").Value = "  +0.91%  "
Set-TextValue "D9" "0.508"
$ws.Range("E9").Value = "  +1.23%  "
Set-TextValue "D10" "7.28"
$ws.Range("E10").Value = "  +3.97%  "
Set-TextValue "D11" "0.151"
$ws.Range("E11").Value = "  +6.79%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("E13").Value = "  +5.22%  "
Set-TextValue "D14" "32.83"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "3.444.88"
$ws.Range("D17").Value = "63.026.57"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "2.963.73"
$ws.Range("E19").Value = "  +1.30%  "
Set-TextValue "D20" "443.50"
$ws.Range("E20").Value = "  +2.72%  "
Set-TextValue "D21" "13.49"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -1.00%  "
Set-TextValue "D23" "7.03"
$ws.Range("E23").Value = "  -0.67%  "
Set-TextValue "D24" "11.25"
$ws.Range("E24").Value = "  +3.46%  "
Set-TextValue "D25" "81.08"
$ws.Range("E25").Value = "  -0.86%  "
Set-TextValue "D26" "2.14"
$ws.Range("E26").Value = "  -2.00%  "
Set-TextValue "D27" "11.77"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.01%  "
Set-TextValue "D29" "7.28"
$ws.Range("E29").Value = "  +5.73%  "
Set-TextValue "D30" "2.21"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +16.74%  "
Set-TextValue "D33" "26.47"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("E34").Value = "  -0.78%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.08%  "
Set-TextValue "D36" "0.990"
$ws.Range("E36").Value = "  -1.87%  "
Set-TextValue "D37" "3.15"
$ws.Range("E37").Value = "  +5.76%  "
$ws.Range("E38").Value = "  -0.60%  "
Set-TextValue "D39" "2.06"
$ws.Range("E39").Value = "  +3.15%  "
Set-TextValue "D40" "49.68"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  -0.57%  "
Set-TextValue "D42" "0.118"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("E43").Value = "  +0.55%  "
Set-TextValue "D44" "38.92"
$ws.Range("E44").Value = "  -7.66%  "
Set-TextValue "D45" "135.56"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D46").Value = "2.693.28"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -2.37%  "
Set-TextValue "D48" "360.62"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.60%  "
Set-TextValue "D51" "22.87"
$ws.Range("E51").Value = "  -3.07%  "
